$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has 5 data rows (2-5) where row 5 is an accidental
# duplicate of row 4 (differing only in column A). This reverts the sheet
# back to 3 distinct data rows (2-4) and fixes up the values that a prior
# "edited via CSV" pass had mangled (path + UUID casing/digits).

# 1. Drop the duplicate last row.
$ws.Rows.Item(5).Delete() | Out-Null

# 2. Row 2: point at the same tileset as the other rows, fix UUID.
$ws.Cells.Item(2, 1).Value = "/tmp/tilesets/TilesetWithDiscreteLOD"
$ws.Cells.Item(2, 3).Value = "b6645aa5-4134-48f3-8cbc-faa0518c21bb"

# 3. Row 3: fix UUID casing/digits (path was already correct).
$ws.Cells.Item(3, 3).Value = "dc61d24d-7426-4490-a48f-06c13be98b85"

# 4. Row 4: point at the same tileset as the other rows, fix UUID.
$ws.Cells.Item(4, 1).Value = "/tmp/tilesets/TilesetWithDiscreteLOD"
$ws.Cells.Item(4, 3).Value = "deda97ad-0912-4524-929c-02beba91c01d"

# 5. Rows 2 and 4 previously had a taller "wrapped text" row height/style
#    applied only because of the stray long path string; now that every
#    row holds the same short path, normalise them back to the standard
#    (non-wrapped) row height/style used by the rest of the sheet (copy
#    the plain formatting from a neighbouring cell in the same row).
$ws.Rows.Item(2).RowHeight = 13.8
$ws.Rows.Item(4).RowHeight = 13.8
$ws.Range("B2").Copy() | Out-Null
$ws.Range("A2").PasteSpecial(-4122) | Out-Null
$ws.Range("B4").Copy() | Out-Null
$ws.Range("A4").PasteSpecial(-4122) | Out-Null

# 6. Restore the original selection (A2) that the CSV round-trip had
#    left pointed at the now-removed row 5.
$ws.Range("A2").Select() | Out-Null
